$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log - Part 2")

# Fill in the newly logged activity rows (18-21) on the Activity Log - Part 2 sheet.
# Row 18: 2020-04-16, 16:10 - 17:15
$ws.Range("B18").Value = 4794
$ws.Range("C18").Value = 43937
$ws.Range("D18").Value = 0.67361111111111116
$ws.Range("E18").Value = 0.71875
$ws.Range("G18").Value = "Finishing screenshots description for RTL and postfit circuits. Completed"

# Row 19: 2020-04-16, 17:15 - 18:45
$ws.Range("B19").Value = 4794
$ws.Range("C19").Value = 43937
$ws.Range("D19").Value = 0.71875
$ws.Range("E19").Value = 0.78125
$ws.Range("G19").Value = "Redoing some screen captures of SLL, SRL and SRA 32 Timing simulation because some info was illegible. "

# Row 20: 2020-04-16, 18:45 - 19:50
$ws.Range("B20").Value = 4794
$ws.Range("C20").Value = 43937
$ws.Range("D20").Value = 0.78125
$ws.Range("E20").Value = 0.82638888888888884
$ws.Range("G20").Value = "Proofreading all documentations and fixing errors. Fixed typos and miscalculations."

# Row 21: 2020-04-16, 20:55 - 22:35
$ws.Range("B21").Value = 4794
$ws.Range("C21").Value = 43937
$ws.Range("D21").Value = 0.87152777777777779
$ws.Range("E21").Value = 0.94097222222222221
$ws.Range("G21").Value = "Working on annotating screenshots. Annotations mostly completed."

# Move the active selection to G16 to match the saved cursor position.
$ws.Activate()
$ws.Range("G16").Select()
